# Update "想去人数" (interested-people count) figures on the 展览, 演出, and
# 全部类型 sheets to match the latest scrape (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1857
$ws1.Range("F6").Value  = 818
$ws1.Range("F16").Value = 4370
$ws1.Range("F19").Value = 478
$ws1.Range("F23").Value = 1793
$ws1.Range("F26").Value = 17
$ws1.Range("F28").Value = 2061
$ws1.Range("F29").Value = 72
$ws1.Range("F32").Value = 146
$ws1.Range("F33").Value = 96
$ws1.Range("F34").Value = 26

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 31

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1857
$ws4.Range("F6").Value  = 818
$ws4.Range("F16").Value = 31
$ws4.Range("F17").Value = 4370
$ws4.Range("F20").Value = 478
$ws4.Range("F24").Value = 1793
$ws4.Range("F27").Value = 17
$ws4.Range("F29").Value = 2061
$ws4.Range("F30").Value = 72
$ws4.Range("F33").Value = 146
$ws4.Range("F34").Value = 96
$ws4.Range("F35").Value = 26
